$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("H2").Value = 1.27
$ws.Range("Q2").Value = 1.41

# --- Row 3 updates ---
$ws.Range("F3").Value = 1.87
$ws.Range("H3").Value = 3.85
$ws.Range("I3").Value = 6.2
$ws.Range("J3").Value = 3.1
$ws.Range("L3").Value = 1.34
$ws.Range("N3").Value = 1.66
$ws.Range("P3").Value = 1.66
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.18
$ws.Range("V3").Value = 1.19
$ws.Range("W3").Value = 1.9

# --- Row 4 updates ---
$ws.Range("F4").Value = 1.7
$ws.Range("K4").Value = 3.95
$ws.Range("N4").Value = 3.3
$ws.Range("P4").Value = 1.74
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.3
$ws.Range("S4").Value = 3.75
$ws.Range("T4").Value = 1.96
$ws.Range("U4").Value = 1.86
$ws.Range("AH4").Value = 23

# --- New row 5 ---
# Use the style of an existing default-styled cell so that forcing text
# number-format to avoid date/time auto-conversion does not leave a
# residual explicit style (cell style index stays the implicit default).
$refStyle = $ws.Range("A4").Style

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "Colombian Primera A"
$ws.Range("A5").Style = $refStyle

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2025-11-12"
$ws.Range("B5").Style = $refStyle

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "22:20:00"
$ws.Range("C5").Style = $refStyle

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "Boyaca Chico"
$ws.Range("D5").Style = $refStyle

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "Millonarios"
$ws.Range("E5").Style = $refStyle

$ws.Range("F5").Value = 4.9
$ws.Range("G5").Value = 6.6
$ws.Range("H5").Value = 1.77
$ws.Range("I5").Value = 1.88
$ws.Range("J5").Value = 3.35
$ws.Range("K5").Value = 3.9
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 2.62
$ws.Range("O5").Value = 1.43
$ws.Range("P5").Value = 1.65
$ws.Range("Q5").Value = 2.24
$ws.Range("R5").Value = 1.24
$ws.Range("S5").Value = 4.3
$ws.Range("T5").Value = 2.1
$ws.Range("U5").Value = 1.76
$ws.Range("V5").Value = 2.12
$ws.Range("W5").Value = 1.18
$ws.Range("X5").Value = 980
$ws.Range("Y5").Value = 980
$ws.Range("Z5").Value = 980
$ws.Range("AA5").Value = 980
$ws.Range("AB5").Value = 980
$ws.Range("AC5").Value = 980
$ws.Range("AD5").Value = 980
$ws.Range("AE5").Value = 980
$ws.Range("AF5").Value = 980
$ws.Range("AG5").Value = 980
$ws.Range("AH5").Value = 980
$ws.Range("AI5").Value = 60
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 980
